# Update the "想去人数" (want-to-go count) figures in column F across the
# three affected sheets, reflecting a newer data scrape.
# 展览 (Exhibitions), 演出 (Performances) and 全部类型 (All types, the
# aggregate sheet) all need their F-column numbers refreshed.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 1594
$wsExhibit.Range("F10").Value = 2657
$wsExhibit.Range("F11").Value = 2657
$wsExhibit.Range("F16").Value = 675
$wsExhibit.Range("F17").Value = 4918
$wsExhibit.Range("F21").Value = 3380
$wsExhibit.Range("F25").Value = 35
$wsExhibit.Range("F26").Value = 2401
$wsExhibit.Range("F27").Value = 55
$wsExhibit.Range("F38").Value = 1360

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F8").Value  = 19
$wsShow.Range("F11").Value = 128
$wsShow.Range("F19").Value = 513

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F14").Value = 1594
$wsAll.Range("F15").Value = 19
$wsAll.Range("F18").Value = 2657
$wsAll.Range("F22").Value = 128
$wsAll.Range("F25").Value = 675
$wsAll.Range("F26").Value = 4918
$wsAll.Range("F29").Value = 3380
$wsAll.Range("F33").Value = 35
$wsAll.Range("F34").Value = 2401
$wsAll.Range("F41").Value = 513
$wsAll.Range("F49").Value = 1360
